$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force Excel to treat the assigned value as literal text even when it
    # looks like a number (e.g. "7.32"), without leaving a visible quote
    # prefix style on the cell once we reset it back to Normal.
    $range.Value = "'" + $value
    $range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "42.512.30"
$ws.Range("E2").Value = "  -0.81%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.530.93"
$ws.Range("E3").Value = "  -0.97%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.05%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "309.02"
$ws.Range("E5").Value = "  -1.97%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "99.68"
$ws.Range("E6").Value = "  +3.69%  "

# Row 7 - XRP
$ws.Range("E7").Value = "  -1.09%  "

# Row 8 - USDC
$ws.Range("E8").Value = "  +0.02%  "

# Row 9 - Cardano
Set-TextValue $ws.Range("D9") "0.528"
$ws.Range("E9").Value = "  -1.68%  "

# Row 10 - Avalanche
Set-TextValue $ws.Range("D10") "35.81"
$ws.Range("E10").Value = "  +1.20%  "

# Row 11 - Dogecoin
Set-TextValue $ws.Range("D11") "0.0803"
$ws.Range("E11").Value = "  -0.93%  "

# Row 12 - Polkadot
Set-TextValue $ws.Range("D12") "7.32"
$ws.Range("E12").Value = "  -1.13%  "

# Row 13 - TRON
$ws.Range("E13").Value = "  +0.15%  "

# Row 14 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D14") "2.914.69"
$ws.Range("E14").Value = "  -1.18%  "

# Row 15 - Chainlink
Set-TextValue $ws.Range("D15") "15.81"
$ws.Range("E15").Value = "  +5.39%  "

# Row 16 - WrappedEther
Set-TextValue $ws.Range("D16") "2.535.03"
$ws.Range("E16").Value = "  -1.32%  "

# Row 17 - Polygon
Set-TextValue $ws.Range("D17") "0.820"
$ws.Range("E17").Value = "  -2.21%  "

# Row 18 - WrappedBTC
Set-TextValue $ws.Range("D18") "42.481.76"
$ws.Range("E18").Value = "  -0.96%  "

# Row 19 - Uniswap
Set-TextValue $ws.Range("D19") "6.82"
$ws.Range("E19").Value = "  +0.13%  "

# Row 20 - ShibaInu
$ws.Range("E20").Value = "  -0.43%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D21") "12.17"
$ws.Range("E21").Value = "  -2.49%  "

# Row 22 - Litecoin
Set-TextValue $ws.Range("D22") "69.07"
$ws.Range("E22").Value = "  +0.00%  "

# Row 23 - BitcoinCash
Set-TextValue $ws.Range("D23") "242.65"
$ws.Range("E23").Value = "  -3.54%  "

# Row 24 - PancakeSwap
$ws.Range("E24").Value = "  -1.99%  "

# Row 25 - ImmutableX
Set-TextValue $ws.Range("D25") "2.04"
$ws.Range("E25").Value = "  -0.84%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  +0.06%  "

# Row 27 - EthereumClassic
Set-TextValue $ws.Range("D27") "25.87"
$ws.Range("E27").Value = "  -3.14%  "

# Row 28 - Toncoin
$ws.Range("E28").Value = "  -4.07%  "

# Row 29 - InjectiveProtocol
Set-TextValue $ws.Range("D29") "39.13"
$ws.Range("E29").Value = "  -1.54%  "

# Row 30 - Cosmos
Set-TextValue $ws.Range("D30") "10.11"
$ws.Range("E30").Value = "  -0.30%  "

# Row 31 - Monero
Set-TextValue $ws.Range("D31") "156.14"
$ws.Range("E31").Value = "  +1.34%  "

# Row 32 - Filecoin
Set-TextValue $ws.Range("D32") "5.73"
$ws.Range("E32").Value = "  -0.78%  "

# Row 33 - ApeXProtocol
Set-TextValue $ws.Range("D33") "2.77"
$ws.Range("E33").Value = "  +13.97%  "

# Row 34 - Hedera
Set-TextValue $ws.Range("D34") "0.0793"
$ws.Range("E34").Value = "  -1.10%  "

# Row 35 - WEMIXToken
$ws.Range("E35").Value = "  -2.90%  "

# Row 36 - ARBITRUM
Set-TextValue $ws.Range("D36") "2.02"
$ws.Range("E36").Value = "  -4.44%  "

# Row 37 - Celestia
Set-TextValue $ws.Range("D37") "18.22"
$ws.Range("E37").Value = "  -4.20%  "

# Row 38 - LidoDAOToken
$ws.Range("E38").Value = "  -6.72%  "

# Row 39 - Kaspa
$ws.Range("E39").Value = "  +0.17%  "

# Row 40 - Stellar
$ws.Range("E40").Value = "  +0.31%  "

# Row 41 - RenderToken
Set-TextValue $ws.Range("D41") "4.30"
$ws.Range("E41").Value = "  +10.31%  "

# Row 42 - EnergySwap
Set-TextValue $ws.Range("D42") "22.04"
$ws.Range("E42").Value = "  -1.81%  "

# Row 43 - FirstDigitalUSD
$ws.Range("E43").Value = "  +0.12%  "

# Row 44 - NEARProtocol
Set-TextValue $ws.Range("D44") "3.31"
$ws.Range("E44").Value = "  +2.30%  "

# Row 45 - VeChain
Set-TextValue $ws.Range("D45") "0.0298"
$ws.Range("E45").Value = "  -1.87%  "

# Row 46 - Maker
Set-TextValue $ws.Range("D46") "1.968.15"
$ws.Range("E46").Value = "  -1.44%  "

# Row 47 - FraxShare
$ws.Range("E47").Value = "  -1.19%  "

# Row 48 - RocketPoolETH
Set-TextValue $ws.Range("D48") "2.768.40"
$ws.Range("E48").Value = "  -1.27%  "

# Row 49 - now becomes SEI (was BitcoinSV)
$ws.Range("B49").Value = "SEI"
$ws.Range("C49").Value = "https://coinranking.com/coin/8nxCqs-uj+sei-sei"
Set-TextValue $ws.Range("D49") "0.864"
$ws.Range("E49").Value = "  +12.53%  "

# Row 50 - now becomes BitcoinSV (was SEI)
$ws.Range("B50").Value = "BitcoinSV"
$ws.Range("C50").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
Set-TextValue $ws.Range("D50") "81.08"
$ws.Range("E50").Value = "  -2.31%  "

# Row 51 - Algorand
Set-TextValue $ws.Range("D51") "0.192"
$ws.Range("E51").Value = "  -0.11%  "
